$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "50.832.59"
$ws.Range("E2").Value = "  -1.15%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.896.16"
$ws.Range("E3").Value = "  -0.91%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "366.03"
$ws.Range("E5").Value = "  +4.42%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "101.70"
$ws.Range("E6").Value = "  -4.07%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -2.71%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.578"
$ws.Range("E9").Value = "  -4.38%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "36.12"
$ws.Range("E10").Value = "  -4.18%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.60%  "

# Row 12 - Dogecoin
Set-TextValue $ws.Range("D12") "0.0825"
$ws.Range("E12").Value = "  -2.62%  "

# Row 13 - now Chainlink (was WrappedliquidstakedEther2.0)
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D13") "18.14"
$ws.Range("E13").Value = "  -4.01%  "

# Row 14 - now WrappedliquidstakedEther2.0 (was Chainlink)
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D14") "3.344.50"
$ws.Range("E14").Value = "  -1.03%  "

# Row 15 - Polkadot
Set-TextValue $ws.Range("D15") "7.32"
$ws.Range("E15").Value = "  -2.76%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.888.51"
$ws.Range("E16").Value = "  -1.17%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.913"
$ws.Range("E17").Value = "  -4.97%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "50.800.95"
$ws.Range("E18").Value = "  -1.13%  "

# Row 19 - ImmutableX
Set-TextValue $ws.Range("D19") "3.19"
$ws.Range("E19").Value = "  -6.41%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "7.11"
$ws.Range("E20").Value = "  -3.75%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "12.71"
$ws.Range("E21").Value = "  -5.09%  "

# Row 22 - ShibaInu
Set-TextValue $ws.Range("D22") "0.0₃0935"
$ws.Range("E22").Value = "  -3.10%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "67.70"
$ws.Range("E23").Value = "  -1.66%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "256.56"
$ws.Range("E24").Value = "  -1.34%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -1.70%  "

# Row 26 - LEO
$ws.Range("E26").Value = "  +2.59%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.01%  "

# Row 28 - Kaspa
Set-TextValue $ws.Range("D28") "0.167"
$ws.Range("E28").Value = "  -3.96%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "25.39"
$ws.Range("E29").Value = "  -3.71%  "

# Row 30 - Filecoin
Set-TextValue $ws.Range("D30") "6.91"
$ws.Range("E30").Value = "  -6.26%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -3.76%  "

# Row 32 - RenderToken
Set-TextValue $ws.Range("D32") "6.10"
$ws.Range("E32").Value = "  +0.69%  "

# Row 33 - Cosmos
Set-TextValue $ws.Range("D33") "9.80"
$ws.Range("E33").Value = "  -3.89%  "

# Row 34 - Toncoin
$ws.Range("E34").Value = "  -3.51%  "

# Row 35 - OKB
Set-TextValue $ws.Range("D35") "50.76"
$ws.Range("E35").Value = "  +0.89%  "

# Row 36 - InjectiveProtocol
Set-TextValue $ws.Range("D36") "33.75"
$ws.Range("E36").Value = "  -5.18%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.44%  "

# Row 38 - VeChain
Set-TextValue $ws.Range("D38") "0.0417"
$ws.Range("E38").Value = "  -2.20%  "

# Row 39 - LidoDAOToken
Set-TextValue $ws.Range("D39") "2.96"
$ws.Range("E39").Value = "  -5.36%  "

# Row 40 - Stacks
Set-TextValue $ws.Range("D40") "2.60"
$ws.Range("E40").Value = "  -1.63%  "

# Row 41 - Celestia
Set-TextValue $ws.Range("D41") "16.85"
$ws.Range("E41").Value = "  -4.30%  "

# Row 42 - ARBITRUM
Set-TextValue $ws.Range("D42") "1.82"
$ws.Range("E42").Value = "  -5.75%  "

# Row 43 - Stellar
$ws.Range("E43").Value = "  -3.43%  "

# Row 44 - Monero
Set-TextValue $ws.Range("D44") "118.57"
$ws.Range("E44").Value = "  -0.88%  "

# Row 45 - EnergySwap
Set-TextValue $ws.Range("D45") "21.72"
$ws.Range("E45").Value = "  -2.22%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  -2.14%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +0.54%  "

# Row 48 - Maker
Set-TextValue $ws.Range("D48") "2.004.39"

# Row 49 - NEARProtocol
Set-TextValue $ws.Range("D49") "3.12"
$ws.Range("E49").Value = "  -5.70%  "

# Row 50 - RocketPoolETH
Set-TextValue $ws.Range("D50") "3.184.05"
$ws.Range("E50").Value = "  -0.64%  "

# Row 51 - TheGraph
$ws.Range("E51").Value = "  -2.00%  "
